# Appends the researcher's log-file rows (35_1_18 / 35_1_19 sessions, plus
# two trailing 2_1_3 rows) to "Sheet 1" of the workbook, directly below the
# existing data (which currently ends at row 30).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$data = @(
    ,(30, 35, 'sml5527', 'poll', 'Q: issFav Q: issOpp Q: candFav Q: candFav Q: undefined Q: undefined', '1490795173254', '35_1_18')
    ,(31, 35, 'sml5527', 'poll', 'Q: issFav Q: issOpp Q: candFav Q: candOpp Q: undefined Q: undefined', '1490795188179', '35_1_18')
    ,(32, 35, 'sml5527', 'poll', 'Q: issFav Q: issOpp Q: candFav Q: undefined Q: candOpp Q: undefined', '1490795244453', '35_1_18')
    ,(33, 35, 'sml5527', 'userAction', 'Work Out', '1490795434798', '35_1_18')
    ,(34, 35, 'sml5527', 'userAction', 'Community Service', '1490795440705', '35_1_18')
    ,(35, 35, 'sml5527', 'userAction', 'Photobomb', '1490796534323', '35_1_18')
    ,(36, 35, 'sml5527', 'userAction', 'Ice Cream Social', '1490796537647', '35_1_18')
    ,(37, 35, 'sml5527', 'userAction', 'Ice Cream Social', '1490796539788', '35_1_18')
    ,(38, 35, 'sml5527', 'userAction', 'Ice Cream Social', '1490796541557', '35_1_18')
    ,(39, 35, 'sml5527', 'userAction', 'Compete in Game Tournament', '1490798333236', '35_1_18')
    ,(40, 35, 'sml5527', 'userAction', 'Ice Cream Social', '1490798335895', '35_1_18')
    ,(41, 35, 'sml5527', 'userAction', 'Ice Cream Social', '1490798338336', '35_1_18')
    ,(42, 35, 'sml5527', 'userAction', 'Ice Cream Social', '1490798340766', '35_1_18')
    ,(43, 35, 'sml5527', 'userAction', 'Ice Cream Social', '1490798345830', '35_1_18')
    ,(44, 35, 'sml5527', 'userAction', 'Ice Cream Social', '1490798353222', '35_1_18')
    ,(45, 35, 'sml5527', 'userAction', 'Compete in Game Tournament', '1490798384022', '35_1_18')
    ,(46, 35, 'sml5527', 'userAction', 'Help Paint Mural', '1490798388102', '35_1_18')
    ,(47, 35, 'sml5527', 'userAction', 'Charity Event', '1490798394634', '35_1_18')
    ,(48, 35, 'sml5527', 'userAction', 'Ice Cream Social', '1490798399735', '35_1_18')
    ,(49, 35, 'sml5527', 'userAction', 'Meet the Robotics Team', '1490798410462', '35_1_18')
    ,(50, 35, 'sml5527', 'endGame', 'Winner: Martian Dog Player Rank:6', '1490798743808', '35_1_18')
    ,(51, 35, 'sml5527', 'poll', 'Q: issFav Q: issOpp Q: candFav Q: candOpp Q: candOpp Q: undefined', '1490799704713', '35_1_19')
    ,(52, 35, 'sml5527', 'poll', 'Q: issFav Q: issOpp Q: candFav Q: candOpp Q: undefined Q: undefined', '1490801937885', '35_1_19')
    ,(53, 35, 'sml5527', 'userAction', 'Help Paint Mural', '1490802115079', '35_1_19')
    ,(54, 35, 'sml5527', 'userAction', 'Create Posters', '1490802123585', '35_1_19')
    ,(55, 35, 'sml5527', 'userAction', 'Compete in Game Tournament', '1490802140752', '35_1_19')
    ,(56, 35, 'sml5527', 'userAction', 'Study in the Library', '1490802146655', '35_1_19')
    ,(57, 35, 'sml5527', 'userAction', 'Community Service', '1490802151791', '35_1_19')
    ,(58, 2, 'aaronjmill32', 'userAction', 'Meet the Debate Team', 'March 30th 2017, 12:41:14 pm', '2_1_3')
    ,(59, 2, 'aaronjmill32', 'userAction', 'Study in the Library', 'March 30th 2017, 12:41:53 pm', '2_1_3')
)

$startRow = 31
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    # The timestamp in column F is a numeric-looking string (epoch millis)
    # that must stay text, matching the rest of the `date` column - prefix
    # with an apostrophe so Excel keeps it as text instead of a number.
    $ws.Cells.Item($r, 6).Value = "'" + $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
